$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88 holds a game entry whose "name" is a number-looking code ("777"),
# so force the cell to text before writing it, then drop back to the
# workbook's Normal style so no stray cell format lingers.
$ws.Cells.Item(88, 1).NumberFormat = "@"
$ws.Cells.Item(88, 1).Value = "777"
$ws.Cells.Item(88, 1).Style = "Normal"

$ws.Cells.Item(88, 2).Value = "Incompleto"
$ws.Cells.Item(88, 3).Value = "Outro"
$ws.Cells.Item(88, 4).Value = "Zerar"
